$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 10000
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -10460
$ws.Range("H35").Value = 10000
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 10000
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 10000
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -10758
$ws.Range("H38").Value = 2637.5881
$ws.Range("I38").Value = 413.16666
$ws.Range("J38").Value = 3850.9092
$ws.Range("K38").Value = 1239.49998
$ws.Range("L38").Value = 11552.7276
$ws.Range("M38").Value = -867.4999800000001
$ws.Range("N38").Value = -12296.7276
$ws.Range("H40").Value = 2991.2
$ws.Range("I40").Value = 6000
$ws.Range("J40").Value = 1701.7142
$ws.Range("K40").Value = 6000
$ws.Range("L40").Value = 1701.7142
$ws.Range("M40").Value = -5825
$ws.Range("N40").Value = -2051.7142
$ws.Range("I76").Value = 3605.2778
$ws.Range("J76").Value = 4045.9092
$ws.Range("K76").Value = 3605.2778
$ws.Range("L76").Value = 4045.9092
$ws.Range("M76").Value = -3290.2778
$ws.Range("N76").Value = -4675.9092
$ws.Range("I79").Value = 3605.2778
$ws.Range("J79").Value = 4045.9092
$ws.Range("K79").Value = 3605.2778
$ws.Range("L79").Value = 4045.9092
$ws.Range("M79").Value = -2513.2778
$ws.Range("N79").Value = -6229.9092
$ws.Range("H86").Value = 104400.3
$ws.Range("I86").Value = 170500.5
$ws.Range("J86").Value = 5250
$ws.Range("K86").Value = 170500.5
$ws.Range("L86").Value = 5250
$ws.Range("M86").Value = -169377.5
$ws.Range("H88").Value = 1008.3182
$ws.Range("I88").Value = 1037.6
$ws.Range("J88").Value = 999.7059
$ws.Range("K88").Value = 1037.6
$ws.Range("L88").Value = 999.7059
$ws.Range("M88").Value = -631.5999999999999
$ws.Range("N88").Value = -1811.7059
$ws.Range("H89").Value = 104400.3
$ws.Range("I89").Value = 170500.5
$ws.Range("J89").Value = 5250
$ws.Range("K89").Value = 852502.5
$ws.Range("L89").Value = 26250
$ws.Range("M89").Value = -846886.5
$ws.Range("H91").Value = 1008.3182
$ws.Range("I91").Value = 1037.6
$ws.Range("J91").Value = 999.7059
$ws.Range("K91").Value = 1037.6
$ws.Range("L91").Value = 999.7059
$ws.Range("M91").Value = 366.4000000000001
$ws.Range("N91").Value = -3807.7059
$ws.Range("H92").Value = 481.08334
$ws.Range("I92").Value = 453.5
$ws.Range("J92").Value = 536.25
$ws.Range("K92").Value = 453.5
$ws.Range("L92").Value = 536.25
$ws.Range("M92").Value = 794.5
$ws.Range("N92").Value = -3032.25
$ws.Range("H137").Value = 3147.122
$ws.Range("I137").Value = 1787.6296
$ws.Range("J137").Value = 5769
$ws.Range("K137").Value = 5362.8888
$ws.Range("L137").Value = 17307
$ws.Range("M137").Value = -2812.8888
$ws.Range("N137").Value = -22407
$ws.Range("H138").Value = 3640.2874
$ws.Range("I138").Value = 2465.15
$ws.Range("J138").Value = 3991.0747
$ws.Range("K138").Value = 7395.450000000001
$ws.Range("L138").Value = 11973.2241
$ws.Range("M138").Value = -2255.450000000001
$ws.Range("N138").Value = -22253.2241

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7506.647
$ws.Range("I32").Value = 6832.5874
$ws.Range("J32").Value = 15999.8
$ws.Range("K32").Value = 6832.5874
$ws.Range("L32").Value = 15999.8
$ws.Range("M32").Value = -6545.5874
$ws.Range("N32").Value = -16573.8
$ws.Range("H132").Value = 4889.5625
$ws.Range("I132").Value = 1897.9565
$ws.Range("J132").Value = 7641.84
$ws.Range("K132").Value = 5693.8695
$ws.Range("L132").Value = 22925.52
$ws.Range("M132").Value = -3163.8695
$ws.Range("N132").Value = -27985.52

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2600
$ws.Range("I86").Value = 2650
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 2650
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -1527
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 2600
$ws.Range("I89").Value = 2650
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 13250
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -7634
$ws.Range("N89").Value = -21232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 1796.6666
$ws.Range("I12").Value = 195
$ws.Range("J12").Value = 5000
$ws.Range("K12").Value = 195
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = -25
$ws.Range("N12").Value = -5340
$ws.Range("H31").Value = 2026.47
$ws.Range("I31").Value = 1428.6418
$ws.Range("J31").Value = 3240.2424
$ws.Range("K31").Value = 1428.6418
$ws.Range("L31").Value = 3240.2424
$ws.Range("M31").Value = -1133.6418
$ws.Range("N31").Value = -3830.2424
$ws.Range("H34").Value = 2026.47
$ws.Range("I34").Value = 1428.6418
$ws.Range("J34").Value = 3240.2424
$ws.Range("K34").Value = 1428.6418
$ws.Range("L34").Value = 3240.2424
$ws.Range("M34").Value = -1226.6418
$ws.Range("N34").Value = -3644.2424
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 8230.742
$ws.Range("I122").Value = 4085
$ws.Range("J122").Value = 11644.883
$ws.Range("K122").Value = 12255
$ws.Range("L122").Value = 34934.649
$ws.Range("M122").Value = -9805
$ws.Range("N122").Value = -39834.649
$ws.Range("H141").Value = 35150.125
$ws.Range("I141").Value = 34148
$ws.Range("J141").Value = 35216.934
$ws.Range("K141").Value = 34148
$ws.Range("L141").Value = 35216.934
$ws.Range("M141").Value = -28968
$ws.Range("N141").Value = -45576.934

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 694.1875
$ws.Range("I122").Value = 381.7619
$ws.Range("J122").Value = 1290.6364
$ws.Range("K122").Value = 3435.8571
$ws.Range("L122").Value = 11615.7276
$ws.Range("M122").Value = -985.8571000000002
$ws.Range("N122").Value = -16515.7276
$ws.Range("H137").Value = 25028778
$ws.Range("I137").Value = 38463816
$ws.Range("J137").Value = 77990.42999999999
$ws.Range("K137").Value = 115391448
$ws.Range("L137").Value = 233971.29
$ws.Range("M137").Value = -115386348
$ws.Range("N137").Value = -244171.29

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 264
$ws.Range("I2").Value = 109.6
$ws.Range("J2").Value = 650
$ws.Range("K2").Value = 109.6
$ws.Range("L2").Value = 650
$ws.Range("M2").Value = 3.400000000000006
$ws.Range("N2").Value = -876
$ws.Range("H70").Value = 5277.282
$ws.Range("I70").Value = 4828
$ws.Range("J70").Value = 5501.923
$ws.Range("K70").Value = 4828
$ws.Range("L70").Value = 5501.923
$ws.Range("M70").Value = -4558
$ws.Range("N70").Value = -6041.923
$ws.Range("H73").Value = 5277.282
$ws.Range("I73").Value = 4828
$ws.Range("J73").Value = 5501.923
$ws.Range("K73").Value = 4828
$ws.Range("L73").Value = 5501.923
$ws.Range("M73").Value = -3892
$ws.Range("N73").Value = -7373.923
$ws.Range("H97").Value = 1931.1765
$ws.Range("I97").Value = 1693.7693
$ws.Range("J97").Value = 2702.75
$ws.Range("K97").Value = 1693.7693
$ws.Range("L97").Value = 2702.75
$ws.Range("M97").Value = -1197.7693
$ws.Range("N97").Value = -3694.75
$ws.Range("H122").Value = 5288.125
$ws.Range("I122").Value = 7801
$ws.Range("J122").Value = 3333.6667
$ws.Range("K122").Value = 23403
$ws.Range("L122").Value = 10001.0001
$ws.Range("M122").Value = -20953

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1885.8572
$ws.Range("I16").Value = 1740.2
$ws.Range("J16").Value = 2250
$ws.Range("K16").Value = 1740.2
$ws.Range("L16").Value = 2250
$ws.Range("M16").Value = -1570.2
$ws.Range("N16").Value = -2590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 4500
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 9000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -11122
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 4500
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 45000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -55608
$ws.Range("H132").Value = 3444.5454
$ws.Range("I132").Value = 2703.9565
$ws.Range("J132").Value = 5147.9
$ws.Range("K132").Value = 8111.869499999999
$ws.Range("L132").Value = 15443.7
$ws.Range("M132").Value = -5581.869499999999
$ws.Range("N132").Value = -20503.7
